$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Pong is one of the earliest arcade video games. It is a table tennis sports game featuring simple two-dimensional graphics. The game was originally manufactured by Atari. which released it in 1972. Allan Alcorn created Pong as a training exercise assigned to him by Atari co-founder Nolan Bushnell. Bushnell based the idea on an electronic ping-pong game included in the Magnavox Odyssey. which later resulted in a lawsuit against Atari. Surprised by the quality of Alcorn''s work. Bushnell and Atari co-founder Ted Dabney decided to manufacture the game.'
$ws.Range("B3").Value = 'Zork is one of the earliest interactive fiction computer games. with roots drawn from the original genre game Colossal Cave Adventure. The first version of Zork was written between 1977 and 1979 using the MDL programming language on a DEC PDP-10 computer. The authors—Tim Anderson. Marc Blank. Bruce Daniels. and Dave Lebling—were members of the MIT Dynamic Modelling Group.'
$ws.Range("B4").Value = 'Space Invaders is an arcade video game created by Tomohiro Nishikado and released in 1978. It was originally manufactured and sold by Taito in Japan. and was later licensed for production in the United States by the Midway division of Bally. Space Invaders is one of the earliest shooting games and the aim is to defeat waves of aliens with a laser cannon to earn as many points as possible. In designing the game. Nishikado drew inspiration from popular media: Breakout. The War of the Worlds. and Star Wars. To complete it. he had to design custom hardware and development tools.'
$ws.Range("B5").Value = 'Asteroids is an arcade space shooter released in November 1979 by Atari. Inc. and designed by Lyle Rains. Ed Logg. and Dominic Walsh. The player controls a spaceship in an asteroid field which is periodically traversed by flying saucers. The object of the game is to shoot and destroy asteroids and saucers while not colliding with either or being hit by the saucers'' counter-fire. The game becomes harder as the number of asteroids increases.'
$ws.Range("B6").Value = 'Adventure is a video game for the Atari 2600 video game console. released in late 1979–1980. In the game. the player controls a square avatar whose quest is to explore an open-ended environment to find a magical chalice and return it to the golden castle. The game world is populated by roaming enemies: three dragons that can eat the avatar and a bat that randomly steals and hides items around the game world. Adventure introduced a number of innovative game elements to console games. including a playing area that spanned several different screens and enemies that continued to move even when not displayed on the screen.'
$ws.Range("B7").Value = 'Pac-Man. stylized as PAC-MAN. is an arcade game developed by Namco and first released in Japan in May 1980.It was created by Japanese video game designer Toru Iwatani. It was licensed for distribution in the United States by Midway Games and released in October 1980. Immensely popular from its original release to the present day. Pac-Man is considered one of the classics of the medium. and an icon of 1980s popular culture. Upon its release. the game—and. subsequently. Pac-Man derivatives—became a social phenomenon that yielded high sales of merchandise and inspired a legacy in other media. such as the Pac-Man animated television series and the top-ten Buckner and Garcia hit single "Pac-Man Fever". Pac-Man was popular in the 1980s and 1990s and is still played in the 2010s.'
$ws.Range("B8").Value = 'Defender is an arcade video game developed and released by Williams Electronics in February 1981. A horizontally scrolling shoot ''em up. the game is set on an unnamed planet where the player must defeat waves of invading aliens while protecting astronauts. Development was led by Eugene Jarvis. a pinball programmer at Williams; Defender was Jarvis'' first video game project and drew inspiration from Space Invaders and Asteroids.'
$ws.Range("B9").Value = 'Donkey Kong is an arcade game released by Nintendo in 1981. An early example of the platform game genre. the gameplay focuses on maneuvering the main character across a series of platforms while dodging and jumping over obstacles. In the game. Mario (originally named Mr. Video and then Jumpman) must rescue a damsel in distress named Pauline (originally named Lady). from a giant ape named Donkey Kong. The hero and ape later became two of Nintendo''s most popular and recognizable characters. Donkey Kong is one of the most important titles from the golden age of arcade video games. and is one of the most popular arcade games of all time.'
$ws.Range("B10").Value = 'Galaga is a Japanese shoot-''em-up arcade game developed and published by Namco Japan and by Midway in North America in 1981. It is the sequel to 1979''s Galaxian. The gameplay of Galaga puts the player in control of a spacecraft which is situated at the bottom of the screen. with enemy aliens arriving in formation at the beginning of a stage. either trying to destroy. collide with. or capture the spaceship. with the player progressing every time alien forces are vanquished.'
$ws.Range("B12").Value = 'Joust is an arcade game developed by Williams Electronics and released in 1982. While not the first game to feature two-player cooperative play. Joust was more successful than its predecessors and popularized the concept. The player uses a button and joystick to control a knight riding a flying ostrich. The object is to progress through levels by defeating groups of enemy knights riding buzzards.'
$ws.Range("B13").Value = 'Ms. Pac-Man is an arcade video game from the Golden Age. It was produced by Illinois-based Midway Manufacturing corporation. the North American publisher of Pac-Man. Ms. Pac-Man was released in North America in January 1982. and is one of the most popular arcade video games of all time. This popularity led to its adoption as an official title by Namco. the creator of Pac-Man. which was released in the United States in late 1980. Ms. Pac-Man introduced a female protagonist. new maze designs. and several other improved gameplay changes over the original Pac-Man. Ms. Pac-Man became the most successful American-produced arcade game of 1981. selling 115.000 arcade cabinets.'
$ws.Range("B15").Value = 'Punch-Out!! is a boxing arcade game by Nintendo. originally released late 1983.[2][8] It was the first in a series of successful Punch-Out!! games. producing an arcade sequel known as Super Punch-Out!!. a spin-off of the series titled Arm Wrestling. a highly popular version for the NES originally known as Mike Tyson''s Punch-Out!!. and Super Punch-Out!! for the SNES.'
$ws.Range("B16").Value = 'Star Wars is an arcade game produced by Atari Inc. and released in 1983. The game is a first person space combat game. simulating the attack on the Death Star from the 1977 film Star Wars. The game is composed of 3D color vector graphics. This game was developed during the Golden Age of Arcade Games and was voted the fourth-best arcade game of all time by the members of Killer List of Videogames.'
$ws.Range("B17").Value = 'Elite is a space trading video game. written and developed by David Braben and Ian Bell and originally published by Acornsoft for the BBC Micro and Acorn Electron computers in September 1984. Elite''s open-ended game model. and revolutionary 3D graphics led to it being ported to virtually every contemporary home computer system. and earned it a place as a classic and a genre maker in gaming history. The game''s title derives from one of the player''s goals of raising their combat rating to the exalted heights of "Elite".'
$ws.Range("B18").Value = 'Tetris (Russian: Те́трис. pronounced [ˈtɛtrʲɪs]) is a tile-matching puzzle video game. originally designed and programmed by Russian game designer Alexey Pajitnov. It was released on June 6. 1984. while he was working for the Dorodnitsyn Computing Centre of the Academy of Science of the Soviet Union in Moscow. He derived its name from the Greek numerical prefix tetra- (all of the game''s pieces contain four segments) and tennis. Pajitnov''s favorite sport.'
$ws.Range("B19").Value = 'Gauntlet is a fantasy-themed hack and slash 1985 arcade game by Atari Games.Released in October 1985. Atari ultimately sold a total of 7.848 Gauntlet video game arcade cabinets. It is noted as being one of the first multi-player dungeon crawl arcade games.The core design of Gauntlet comes from Dandy. a 1983 Atari 8-bit family title. which resulted in a lawsuit.'
$ws.Range("B20").Value = 'Ghosts ''n Goblins . is a 1985 side-scrolling platform video game developed by Capcom for video arcades and has since been released on several other platforms. It is the first game in the Ghosts ''n Goblins franchise. It was directed by Tokuro Fujiwara.'
$ws.Range("B21").Value = 'Super Mario Bros. is a platform video game developed and published by Nintendo for the Nintendo Entertainment System home console. Released as a successor to the 1983 game Mario Bros.. Super Mario Bros. was released in Japan and North America in 1985. and in Europe and Australia two years later. In Super Mario Bros.. the player controls Mario and his brother. Luigi. as they travel through the Mushroom Kingdom in order to rescue Princess Toadstool from the antagonist. Bowser.'
$ws.Range("B22").Value = 'Final Fantasy VII is a role-playing video game developed by Square for the PlayStation console. Released in 1997. it is the seventh main installment in the Final Fantasy series. Published in Japan by Square. it was released in the West by Sony Computer Entertainment. becoming the first in the main series to be released in Europe. The game''s story follows Cloud Strife. a mercenary who joins an eco-terrorist organization to stop the world-controlling megacorporation. Shinra. from using the planet''s life essence as an energy source. Events send Cloud and his allies in pursuit of Sephiroth. a superhuman intent on destroying their planet. During the journey. Cloud builds close friendships with his party members. including Aerith Gainsborough. who holds the secret to saving their world.'
$ws.Range("B23").Value = 'Final Fantasy IX is a role-playing video game developed and published by Square for the PlayStation video game console. Originally released in 2000. it is the ninth title in the Final Fantasy series and last to debut on the original PlayStation console. In 2010. it was re-released worldwide as a PSOne Classic on the PlayStation Store compatible with PlayStation 3 and PlayStation Portable with PlayStation Vita support coming in 2012. The game introduced new features to the series like the "Active Time Event". "Mognet". and a unique equipment and skill system.'
$ws.Range("B25").Value = 'Donkey Kong Country 2: Diddy''s Kong Quest. plattformsspel som utvecklades av Rare till Super Nintendo. och SNES-versionen släpptes den 21 november 1995 i Japan. Spelet släpptes även till Game Boy Advance 2004.'
$ws.Range("B26").Value = ' Space Invaders Part II is a fixed shooter arcade game which was released by Taito Corporation in 1979; it utilizes one Intel 8080 microprocessor (running at: 2 MHz). with discrete components and a Texas Instruments SN-76477 for audio. It was later licensed to Midway Games for US manufacture and distribution in 1980 under the name Space Invaders Deluxe - however. the title screen on this version still says "Space Invaders Part II". when it is displayed during the game''s attract sequence. This version also reduced the point value of the new "flashing" spaceships. from 500 points to 200.'
